$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Crime Complaints table updates (rows 15-30) ---
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 2
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = -75
$ws.Range("H15").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("L15").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 9
$ws.Range("K16").Value = 125
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -64
$ws.Range("N16").Value = -92.436974789916
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 14
$ws.Range("K17").Value = -28.571428571428
$ws.Range("L17").Value = 11.111111111111
$ws.Range("M17").Value = 11.111111111111
$ws.Range("N17").Value = -23.076923076923
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 18.181818181818
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = 52.941176470588
$ws.Range("L18").Value = 52.941176470588
$ws.Range("M18").Value = 73.333333333333
$ws.Range("N18").Value = -88.288288288288
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 38
$ws.Range("I19").Value = 57
$ws.Range("J19").Value = 69
$ws.Range("K19").Value = -17.391304347826
$ws.Range("L19").Value = 42.5
$ws.Range("M19").Value = -6.55737704918
$ws.Range("N19").Value = -55.11811023622
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 83.333333333333
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = 90.90909090909
$ws.Range("L20").Value = 250
$ws.Range("M20").Value = 23.529411764705
$ws.Range("N20").Value = -96.209386281588
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -10.526315789473
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = 10.144927536231
$ws.Range("I21").Value = 124
$ws.Range("J21").Value = 119
$ws.Range("K21").Value = 4.201680672268
$ws.Range("L21").Value = 51.219512195122
$ws.Range("M21").Value = -2.362204724409
$ws.Range("N21").Value = -88.042430086788
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -28.571428571428
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 25
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = -50
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -12.903225806451
$ws.Range("F24").Value = 136
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -4.225352112676
$ws.Range("I24").Value = 231
$ws.Range("J24").Value = 240
$ws.Range("K24").Value = -3.75
$ws.Range("L24").Value = 4.054054054054
$ws.Range("M24").Value = 133.333333333333
$ws.Range("F15").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 23.076923076923
$ws.Range("I25").Value = 29
$ws.Range("J25").Value = 19
$ws.Range("K25").Value = 52.631578947368
$ws.Range("L25").Value = 45
$ws.Range("M25").Value = -9.375
$ws.Range("F15").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 2
$ws.Range("H15").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = -60
$ws.Range("H15").Copy()
$ws.Range("L26").PasteSpecial(-4122)
$ws.Range("L26").Value = 100
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = -20
$ws.Range("D30").Value = 2
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -83.333333333333
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -83.333333333333
$ws.Range("L30").Value = -50

$excel.CutCopyMode = $false
